$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.387.48'
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").Value = '3.325.20'
$ws.Range("E3").Value = '  -4.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.86'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.98%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.80%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '3.325.69'
$ws.Range("E9").Value = '  -4.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.127'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.84'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("E12").Value = '  -0.77%  '

$ws.Range("D13").Value = '3.904.75'
$ws.Range("E13").Value = '  -4.16%  '

$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.34'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.89%  '

$ws.Range("D16").Value = '65.423.31'
$ws.Range("E16").Value = '  -0.83%  '

$ws.Range("E17").Value = '  -1.71%  '

$ws.Range("D18").Value = '3.322.83'
$ws.Range("E18").Value = '  -4.56%  '

$ws.Range("E19").Value = '  -3.22%  '

$ws.Range("E20").Value = '  -3.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.04'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.38'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.06'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.515'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.17%  '

$ws.Range("E26").Value = '  -3.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.79%  '

$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("E30").Value = '  -1.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.56'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.81'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.00'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("E37").Value = '  -3.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.842'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.41'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.39%  '

$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("D41").Value = '2.696.91'
$ws.Range("E41").Value = '  -4.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.20'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.69%  '

$ws.Range("E44").Value = '  -4.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0664'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '331.65'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.78'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0276'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.72%  '

$ws.Range("E50").Value = '  +1.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.01%  '
